$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ficha técnica")

# Remove the "DIMENSIÓN" / "Accesibilidad" row (row 3) - rows below shift up
$ws.Rows.Item(3).Delete()

# Append new metadata rows at the end of the table (now rows 7 and 8)
$ws.Cells.Item(7, 1).Value2 = "TIPOIND"
$ws.Cells.Item(7, 2).Value2 = "Resultados"

$ws.Cells.Item(8, 1).Value2 = "CITA"
$ws.Cells.Item(8, 2).Value2 = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE"
